# Generate Report for Handback
# Update the timestamp strings recorded in the handback status workbook.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first data row.
# (This timestamp is shared with de-de's "Correspond Handoff Datetime" below,
#  since both describe the same handoff event for the same file.)
$wsOverview.Range("G2").Value = "2016-08-23 01:02:20"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the first data row.
$wsZhCn.Range("H2").Value = "2016-08-23 01:02:14"
$wsZhCn.Range("K2").Value = "2016-08-23 01:02:31"

# de-de sheet: "Correspond Handoff Datetime" (mirrors Overview!G2) and
# "Correspond Handback DateTime" for the first data row.
$wsDeDe.Range("H2").Value = "2016-08-23 01:02:20"
$wsDeDe.Range("K2").Value = "2016-08-23 01:02:39"
